# Bug fix + feedback update: refresh the feedback export with the two new
# Runrunit tickets (FBMDS, Litero) and drop the rows that scrolled out of
# the reporting window so the sheet matches the fresh pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New / updated feedback rows (id, titulo, estado, Quadro, tags, ocorridos, data, cliente)
$rows = @(
    @(156244, "FBMDS, 29/07/2024",      "backlog", "Acompanhamento de clientes", "[]", "Planejamento foi reprovado",                              "2024-07-29", "FBMDS"),
    @(156244, "FBMDS, 29/07/2024",      "backlog", "Acompanhamento de clientes", "[]", "Houve atraso nas entregas, isso prejudicou os clientes",  "2024-07-29", "FBMDS"),
    @(156243, "Litero, 05/08/2024",     "backlog", "Acompanhamento de clientes", "[]", "Cliente pediu proposta",                                  "2024-08-05", "Litero"),
    @(156243, "Litero, 05/08/2024",     "backlog", "Acompanhamento de clientes", "[]", "Feedback positivo",                                       "2024-08-05", "Litero"),
    @(154950, "Mart Minas, 12/08/2024", "backlog", "Acompanhamento de clientes", "[]", "Houve atraso nas entregas, isso prejudicou os clientes",  "2024-08-12", "Mart Minas"),
    @(154950, "Mart Minas, 12/08/2024", "backlog", "Acompanhamento de clientes", "[]", "Cliente solicitou ajustes ou refação",                    "2024-08-12", "Mart Minas"),
    @(153488, "Mart Minas, 07/08/2024", "backlog", "Acompanhamento de clientes", "[]", "Cliente solicitou ajustes ou refação",                    "2024-08-07", "Mart Minas"),
    @(153488, "Mart Minas, 07/08/2024", "backlog", "Acompanhamento de clientes", "[]", "Entregas feitas conforme planejado",                      "2024-08-07", "Mart Minas"),
    @(149896, "Mart Minas, 15/07/2024", "backlog", "Acompanhamento de clientes", "[]", "Metas não atingidas",                                     "2024-07-15", "Mart Minas"),
    @(149896, "Mart Minas, 15/07/2024", "backlog", "Acompanhamento de clientes", "[]", "Feedback positivo",                                       "2024-07-15", "Mart Minas")
)

# Column G holds dates formatted as plain text (yyyy-mm-dd); force text
# formatting first so Excel doesn't silently coerce them into date serials.
$ws.Range("G2:G11").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value2 = $data[0]
    $ws.Cells.Item($r, 2).Value2 = $data[1]
    $ws.Cells.Item($r, 3).Value2 = $data[2]
    $ws.Cells.Item($r, 4).Value2 = $data[3]
    $ws.Cells.Item($r, 5).Value2 = $data[4]
    $ws.Cells.Item($r, 6).Value2 = $data[5]
    $ws.Cells.Item($r, 7).Value2 = $data[6]
    $ws.Cells.Item($r, 8).Value2 = $data[7]
}

# The older feedback rows (previously rows 12-17) are outside the window now;
# drop them so the sheet ends at row 11 (dimension A1:H11).
$ws.Rows("12:17").Delete()
